$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header labels for new columns D (completed) and E (remaining).
# Shared-string table order matters: "remaining" must land before
# "completed" so the new unique strings append in that sequence.
$ws.Range("E2").Value = "remaining"
$ws.Range("D2").Value = "completed"

# Per-row completed(D)/remaining(E) hour splits
$data = @(
    @(4, 0, 3),
    @(5, 0, 3),
    @(6, 2, 1),
    @(7, 3, 0),
    @(11, 2, $null),
    @(14, 1, $null),
    @(15, 2, $null),
    @(16, $null, 0.5),
    @(18, $null, 0.5),
    @(19, $null, 0.5),
    @(20, $null, 0.5),
    @(22, $null, 0.5),
    @(23, $null, 0.5),
    @(24, $null, 0.5),
    @(27, 3, $null),
    @(28, 2, 2),
    @(31, 1, $null),
    @(32, 2, 1),
    @(35, $null, 2),
    @(36, $null, 2),
    @(37, $null, 0.5),
    @(39, $null, 3),
    @(40, $null, 3),
    @(41, $null, 0.5),
    @(44, 0, 1),
    @(45, 0, 1),
    @(46, $null, 0.5),
    @(48, 2, $null),
    @(49, 2, $null),
    @(50, $null, 0.5),
    @(52, 1, 0),
    @(53, 1, 1),
    @(54, $null, 0.5),
    @(56, 1, 0),
    @(57, 1, 1),
    @(58, $null, 0.5),
    @(60, 0, 1),
    @(61, 0, 2),
    @(62, 0, 0.5),
    @(64, 0, 3.5),
    @(65, 0, 4),
    @(66, 0, 0.5),
    @(68, 0, 2),
    @(69, 0, 2),
    @(70, 0, 0.5),
    @(74, 1, 1),
    @(75, 0, 0.5),
    @(77, 2, 0),
    @(78, 0, 0.5),
    @(80, 2, 0),
    @(81, 0, 0.5),
    @(83, 1, 1),
    @(84, 0, 0.5),
    @(86, 0, 2),
    @(87, 0, 0.5),
    @(88, 2, 0),
    @(89, 2, 0),
    @(92, 0, 5),
    @(93, 2, 0)
)

foreach ($row in $data) {
    $r = $row[0]
    $d = $row[1]
    $e = $row[2]
    if ($null -ne $d) {
        $ws.Cells.Item($r, 4).Value = $d
    }
    if ($null -ne $e) {
        $ws.Cells.Item($r, 5).Value = $e
    }
}

# Totals row: SUM formulas for completed/remaining columns
$ws.Range("D95").Formula = "=SUM(D4:D94)"
$ws.Range("E95").Formula = "=SUM(E4:E94)"

# Column D width (auto matches content after the new data)
$ws.Columns.Item(4).ColumnWidth = 10.5546875

# Selection / view state from the diff
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("M6").Select()
